$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Staff Vitals"
#   - remove the "First Name" column (old B)
#   - insert a new "Staff Vitals - FIRSTNAME" column at (new) K
#   - header row loses its bold/bordered/centered formatting
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Staff Vitals")
$ws1.Columns.Item(2).Delete()
$ws1.Columns.Item(11).Insert()
$ws1.Cells.Item(1, 11).Value = "Staff Vitals - FIRSTNAME"
$ws1.Range("A1:AA1").ClearFormats()

# ---------------------------------------------------------------------------
# Sheet 2: "Staff Attributes"
#   - content is unchanged, only the header formatting is cleared
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Staff Attributes")
$ws2.Range("A1:W1").ClearFormats()

# ---------------------------------------------------------------------------
# Sheet 3: "Staff Style"
#   - remove "Grit & Grind Proficiency" (old C), "Pace &Space Proficiency"
#     (old D) and "Seven Seconds Proficiency" (old F)
#   - insert "Staff Style - GRIT_&_GRIND_PROFICIENCY" (new E),
#     "Staff Style - PACE_&_SPACE_PROFICIENCY" (new I) and
#     "Staff Style - SEVEN_SECONDS_PROFICIENCY" (new M)
#   - header row loses its bold/bordered/centered formatting
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Staff Style")
$ws3.Columns.Item(6).Delete()
$ws3.Columns.Item(4).Delete()
$ws3.Columns.Item(3).Delete()

$ws3.Columns.Item(5).Insert()
$ws3.Cells.Item(1, 5).Value = "Staff Style - GRIT_&_GRIND_PROFICIENCY"

$ws3.Columns.Item(9).Insert()
$ws3.Cells.Item(1, 9).Value = "Staff Style - PACE_&_SPACE_PROFICIENCY"

$ws3.Columns.Item(13).Insert()
$ws3.Cells.Item(1, 13).Value = "Staff Style - SEVEN_SECONDS_PROFICIENCY"

$ws3.Range("A1:Q1").ClearFormats()

# ---------------------------------------------------------------------------
# Sheet 4: "Staff Coaching"
#   - content is unchanged, only the header formatting is cleared
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Staff Coaching")
$ws4.Range("A1:B1").ClearFormats()
